$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: correct team member name to "Sourabh Sing" (with leading line break,
# matching the existing shared string used elsewhere), wrap the text and
# increase the row height to fit two lines.
$ws.Range("D22").Value = "`nSourabh Sing"
$ws.Range("D22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 28.8

# New row 26: "Insert statement for Module" task assigned to Ritika Kumar.
$ws.Range("A26").Value = "Insert statement for Module"
$ws.Range("D26").Value = "Ritika Kumar"

# Turn on AutoFilter for the full table range, which also registers the
# hidden _xlnm._FilterDatabase defined name scoped to this sheet.
[void]$ws.Range("A1:E26").AutoFilter()
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$26")
$n.Visible = $false

# Restore the active cell selection to match the saved view.
[void]$ws.Range("B18").Select()
